# The AP210 (ISO 10303-410) entries were moved out of this CR's tracking
# table (they now belong to CR_designAPs_1), so the three rows covering
# N9374 / N9375 / N9376 (the ed4 ap210_electronic_assembly_interconnect_
# and_packaging_design Document / ARM EXPRESS / MIM EXPRESS items) are
# removed from the "WG NB" worksheet. Deleting the entire rows (rather
# than just clearing cell contents) shifts every following row up by
# three, which also drops the now-unused shared strings and shrinks the
# sheet's used range from A1:M39 to A1:M36.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:4").Delete()

# Leave the active selection where the author's commit left it.
$ws.Range("E8").Select()
